# Generate Report for Handoff
# Regenerates the localization-status report: the handoff markdown file was
# re-rolled with a new GUID-named commit (eff3bdf3... -> e5b5e060...), and the
# handoff timestamps/target filenames were refreshed accordingly.

$wb = $excel.ActiveWorkbook

$oldGuid = "eff3bdf3-64a9-482d-ad8a-48f3228bca65"
$newGuid = "e5b5e060-5285-4db9-8e1f-6904a89979b8"

$oldMd   = "$oldGuid.md"
$newMd   = "$newGuid.md"
$oldMdPath = "e2e\$oldGuid.md"
$newMdPath = "e2e\$newGuid.md"

$oldZhXlf = "$oldGuid.858f551ad92d2debb66ff536e8be7e4e3cae3b6e.zh-cn.xlf"
$newZhXlf = "$newGuid.8c0aaa83b6b191bcbcd38718b53fc24f95f2532b.zh-cn.xlf"
$oldDeXlf = "$oldGuid.858f551ad92d2debb66ff536e8be7e4e3cae3b6e.de-de.xlf"
$newDeXlf = "$newGuid.8c0aaa83b6b191bcbcd38718b53fc24f95f2532b.de-de.xlf"

$oldHoDate = "2016-09-04 11:03:08"
$newHoDate = "2016-09-04 11:03:33"
$oldZhDate = "2016-09-04 11:02:59"
$newZhDate = "2016-09-04 11:03:28"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f2d881ebd78ca605c06375e6bf55149c60af925/e2e/$oldGuid.md"

# Original "HyperLink" cell-style colour (RGB 6495ED), as an OLE BGR integer,
# so the re-created hyperlinks keep the same look instead of falling back to
# the generic theme hyperlink colour.
$hyperlinkColor = 15570276

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", $newMdPath)
$wsOverview.Range("B2").Font.Color = $hyperlinkColor

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, "", "", $newMd)
$wsZhCn.Range("A2").Font.Color = $hyperlinkColor

# --- de-de sheet ------------------------------------------------------------
$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newHoDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, "", "", $newMd)
$wsDeDe.Range("A2").Font.Color = $hyperlinkColor
